# Applies the "Elimna EC anteriores y se agregan nuevos, se modifica base de datos" edit:
#  - Removes the previous set of account-statement detail rows and writes a new set
#    (8 worker/period records instead of 7), including a brand-new trailing record.
#  - Updates the totals block (VALOR MORA, Cant. Trabajadores, Cant. Periodos).
#  - Keeps the signature block at the bottom, now shifted one row down to make room
#    for the extra detail row.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# ------------------------------------------------------------------
# 1. Make room for one extra detail row in the table (old last row was 22,
#    new last row is 23). Inserting a row at 23 pushes the signature block
#    (rows 27-28) down to rows 28-29, exactly like the target layout.
# ------------------------------------------------------------------
$ws.Rows.Item(23).Insert(-4121, 0)   # xlShiftDown, no format copied from insert itself

# Propagate the correct per-row formatting:
#  - row 22 currently still carries the special "last row" border/style -> move it to row 23
#  - row 21 carries the normal interior-row border/style -> copy it onto row 22
$ws.Range("B22:J22").Copy($ws.Range("B23:J23"))
$ws.Range("B21:J21").Copy($ws.Range("B22:J22"))

# ------------------------------------------------------------------
# 2. Rewrite the detail table (rows 16-23) with the new data set.
# ------------------------------------------------------------------
$tipoDoc = "CC"

$ws.Range("B16").Value = $tipoDoc
$ws.Range("C16").Value = "1002412781"
$ws.Range("D16").Value = "JEAN CARLOS SARMIENTO AMARANTO"
$ws.Range("E16").Value = "2107"
$ws.Range("F16").Value = 36341
$ws.Range("G16").Value = 908526

$ws.Range("B17").Value = $tipoDoc
$ws.Range("C17").Value = "1002412781"
$ws.Range("D17").Value = "JEAN CARLOS SARMIENTO AMARANTO"
$ws.Range("E17").Value = "2106"
$ws.Range("F17").Value = 36341
$ws.Range("G17").Value = 908526

$ws.Range("B18").Value = $tipoDoc
$ws.Range("C18").Value = "1002412778"
$ws.Range("D18").Value = "VRIGILIO SARMIENTO AMARANTO"
$ws.Range("E18").Value = "2107"
$ws.Range("F18").Value = 36341
$ws.Range("G18").Value = 908526

$ws.Range("B19").Value = $tipoDoc
$ws.Range("C19").Value = "1002412778"
$ws.Range("D19").Value = "VRIGILIO SARMIENTO AMARANTO"
$ws.Range("E19").Value = "2106"
$ws.Range("F19").Value = 36341
$ws.Range("G19").Value = 908526

$ws.Range("B20").Value = $tipoDoc
$ws.Range("C20").Value = "8602824"
$ws.Range("D20").Value = "FABIAN ALMANZA JIMENEZ"
$ws.Range("E20").Value = "2209"
$ws.Range("F20").Value = 40000
$ws.Range("G20").Value = 1000000

$ws.Range("B21").Value = $tipoDoc
$ws.Range("C21").Value = "1002412779"
$ws.Range("D21").Value = "EDILBERTO SARMIENTO AMARANTO"
$ws.Range("E21").Value = "2107"
$ws.Range("F21").Value = 36341
$ws.Range("G21").Value = 908526

$ws.Range("B22").Value = $tipoDoc
$ws.Range("C22").Value = "1002412779"
$ws.Range("D22").Value = "EDILBERTO SARMIENTO AMARANTO"
$ws.Range("E22").Value = "2106"
$ws.Range("F22").Value = 36341
$ws.Range("G22").Value = 908526

$ws.Range("B23").Value = $tipoDoc
$ws.Range("C23").Value = "1044921850"
$ws.Range("D23").Value = "IVAN RENE ZAMBRANO GONZALEZ"
$ws.Range("E23").Value = "2106"
$ws.Range("F23").Value = 36341
$ws.Range("G23").Value = 781242

# ------------------------------------------------------------------
# 3. Update the totals block above the table.
# ------------------------------------------------------------------
$ws.Range("E11").Value = 294387   # VALOR MORA total
$ws.Range("C13").Value = 5        # Cant. Trabajadores
$ws.Range("F13").Value = 3        # Cant. Periodos

Write-Host "Edit applied: table rewritten (rows 16-23), totals updated, signature block shifted."
